$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.231.29'
$ws.Range("E2").Value = '  -0.25%  '

# Row 3
$ws.Range("D3").Value = '3.675.01'
$ws.Range("E3").Value = '  -0.36%  '

# Row 4
$ws.Range("E4").Value = '  +0.15%  '

# Row 5
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '678.81'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -0.40%  '

# Row 6
$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '157.59'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  -3.07%  '

# Row 7
$ws.Range("E7").Value = '  -0.04%  '

# Row 8
$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.493'
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  -1.14%  '

# Row 9
$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.146'
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  -1.28%  '

# Row 10
$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '6.93'
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  -5.25%  '

# Row 11
$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.436'
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  -2.11%  '

# Row 12
$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '0.0000231'
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  -3.32%  '

# Row 13
$ws.Range("D13").Value = '4.297.13'
$ws.Range("E13").Value = '  -0.35%  '

# Row 14
$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '32.22'
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  -3.77%  '

# Row 15
$ws.Range("D15").Value = '3.678.07'
$ws.Range("E15").Value = '  -0.31%  '

# Row 16
$ws.Range("D16").Value = '69.227.81'
$ws.Range("E16").Value = '  -0.28%  '

# Row 17
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '0.114'
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  +1.19%  '

# Row 18
$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '16.00'
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  -1.47%  '

# Row 19
$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '6.40'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  -3.42%  '

# Row 20
$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '467.69'
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  -3.29%  '

# Row 21
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '9.96'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  +0.79%  '

# Row 22
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '0.647'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  -2.69%  '

# Row 23
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '79.90'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  -0.48%  '

# Row 24
$ws.Range("D24").Value = '3.824.28'
$ws.Range("E24").Value = '  -0.28%  '

# Row 25
$ws.Range("E25").Value = '  -0.08%  '

# Row 26
$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '0.0000121'
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  -7.26%  '

# Row 27
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '10.90'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -5.25%  '

# Row 28
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '8.99'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  -6.76%  '

# Row 29
$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '2.68'
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -1.71%  '

# Row 30
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '1.73'
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  -6.32%  '

# Row 31
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '6.58'
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  -3.67%  '

# Row 32
$ws.Range("B32").Value = 'Binance-PegBSC-USD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '1.00'
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  +0.13%  '

# Row 33
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '26.90'
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  -0.67%  '

# Row 34
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '1.97'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  -5.63%  '

# Row 35
$ws.Range("B35").Value = 'RenzoRestakedETH'
$ws.Range("C35").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D35").Value = '3.668.31'
$ws.Range("E35").Value = '  +0.32%  '

# Row 36
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '0.162'
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  -3.21%  '

# Row 37
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '8.14'
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  -4.30%  '

# Row 38
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '6.13'
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  -3.03%  '

# Row 39
$ws.Range("E39").Value = '  +0.01%  '

# Row 40
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  -0.05%  '

# Row 41
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '2.22'
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  -1.22%  '

# Row 42
$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '0.0897'
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  -3.88%  '

# Row 43
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '173.67'
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  +8.99%  '

# Row 44
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '0.940'
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  -1.53%  '

# Row 45
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '47.50'
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  -2.03%  '

# Row 46
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '0.000276'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  -4.43%  '

# Row 47
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '2.66'
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  -7.49%  '

# Row 48
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '27.26'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  -7.97%  '

# Row 49
$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '1.27'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  -6.99%  '

# Row 50
$ws.Range("B50").Value = 'SuiNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '1.08'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  -1.90%  '

# Row 51
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '7.76'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  -3.45%  '

